# Apply "Penalty Reward System" (unfinished) edits:
#  - Forecast Comparison: shift Week_Start_Date (col B) forward one week,
#    and overwrite MyForecast (col D) with new (much smaller) values.
#  - Summary: update several derived/summary metrics in col B.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Forecast Comparison ----
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Make sure column B keeps being stored as plain text (dates are text
# strings in this workbook, not real Excel dates) by forcing a text
# number format before assigning the values.
$ws1.Range("B2:B17").NumberFormat = "@"

$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$newForecast = @(1, 1, 1, 2, 2, 1, 1, 1, 1, 1, 1, 2, 1, 2, 2, 1)

for ($i = 0; $i -lt 16; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 2).Value = $newDates[$i]
    $ws1.Cells.Item($row, 4).Value = $newForecast[$i]
}

# ---- Sheet: Summary ----
$ws2 = $wb.Worksheets.Item("Summary")

# Column B on this sheet holds every value as text (even numbers), so
# force text formatting before writing so nothing gets auto-coerced to a
# number or date serial.
$ws2.Range("B2:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"
$ws2.Range("B4").Value = "76"
$ws2.Range("B6").Value = "22"
$ws2.Range("B8").Value = "2376 units"
$ws2.Range("B9").Value = "23"
$ws2.Range("B10").Value = "11"
$ws2.Range("B11").Value = "5"
$ws2.Range("B12").Value = "2"
$ws2.Range("B13").Value = "2025-03-30"
$ws2.Range("B14").Value = "1"
$ws2.Range("B15").Value = "2025-01-19"
